$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '308.95'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '-0.56%'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '37.51'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '-0.05%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.171'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '1.73%'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '1.59%'
$ws.Range('B6').NumberFormat = "@"
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').NumberFormat = "@"
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.911'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '1.60%'
$ws.Range('B7').NumberFormat = "@"
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').NumberFormat = "@"
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '8.247'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '0.54%'
$ws.Range('B8').NumberFormat = "@"
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').NumberFormat = "@"
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.000'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '1.99%'
$ws.Range('B9').NumberFormat = "@"
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').NumberFormat = "@"
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.9368'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '1.79%'
$ws.Range('B10').NumberFormat = "@"
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').NumberFormat = "@"
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1119'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '-6.03%'
$ws.Range('B11').NumberFormat = "@"
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').NumberFormat = "@"
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.1967'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '2.54%'
$ws.Range('B12').NumberFormat = "@"
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').NumberFormat = "@"
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.09022'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '1.66%'
$ws.Range('B13').NumberFormat = "@"
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').NumberFormat = "@"
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.03303'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '-2.59%'
$ws.Range('B14').NumberFormat = "@"
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').NumberFormat = "@"
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.09596'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '-1.11%'
$ws.Range('B15').NumberFormat = "@"
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').NumberFormat = "@"
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.001386'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '0.68%'
$ws.Range('B16').NumberFormat = "@"
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').NumberFormat = "@"
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.006209'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '8.51%'
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').NumberFormat = "@"
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.614'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '1.77%'
$ws.Range('B18').NumberFormat = "@"
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').NumberFormat = "@"
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '4.423'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '1.77%'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.3412'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '0.13%'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.425'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '27.58%'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.1281'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '0.88%'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.2519'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '-2.75%'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04410'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '0.56%'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.001234'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '1.83%'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.004575'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '7.81%'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0001363'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '0.94%'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.02205'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '4.19%'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.05090'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '2.80%'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.007461'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '-2.59%'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1353'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '0.88%'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.008766'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '-11.55%'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.002135'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '3.61%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.008629'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '-10.27%'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00006549'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '-0.17%'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.00000000752'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '0.19%'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.002863'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '-5.91%'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '-40.66%'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '0.19%'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0002004'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '0.19%'
